# Update 2022 (column I) violent-crime counts to reflect data through 2022-12-22.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 7120
$ws.Range("I3").Value = 7330
$ws.Range("I4").Value = 1690
$ws.Range("I5").Value = 695
$ws.Range("I6").Value = 8745
$ws.Range("I7").Value = 25580

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I3").Value = 69
$ws.Range("I7").Value = 299

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 237
$ws.Range("I5").Value = 22
$ws.Range("I6").Value = 241
$ws.Range("I7").Value = 795

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 361
$ws.Range("I7").Value = 967

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 73
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 190
$ws.Range("I7").Value = 597

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I5").Value = 76
$ws.Range("I7").Value = 801
$ws.Range("I8").Value = 1514
$ws.Range("I11").Value = 391
$ws.Range("I15").Value = 294
$ws.Range("I18").Value = 201
$ws.Range("I19").Value = 715
$ws.Range("I20").Value = 633
$ws.Range("I33").Value = 1123
$ws.Range("I37").Value = 795
$ws.Range("I41").Value = 111
$ws.Range("I42").Value = 974
$ws.Range("I44").Value = 195
$ws.Range("I47").Value = 185
$ws.Range("I48").Value = 325
$ws.Range("I51").Value = 295
$ws.Range("I52").Value = 578
$ws.Range("I53").Value = 293
$ws.Range("I55").Value = 304
$ws.Range("I60").Value = 148
$ws.Range("I63").Value = 79
$ws.Range("I65").Value = 597
$ws.Range("I67").Value = 967
$ws.Range("I71").Value = 75
$ws.Range("I73").Value = 231
$ws.Range("I74").Value = 40
$ws.Range("I78").Value = 338
$ws.Range("I84").Value = 223
$ws.Range("I85").Value = 1136
$ws.Range("I89").Value = 299
$ws.Range("I90").Value = 334
$ws.Range("I91").Value = 271
$ws.Range("I93").Value = 146
$ws.Range("I95").Value = 395
$ws.Range("I98").Value = 185
$ws.Range("I101").Value = 25580

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 141
$ws.Range("I6").Value = 83
$ws.Range("I7").Value = 395

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 255
$ws.Range("I7").Value = 1123

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 517
$ws.Range("I6").Value = 423

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 207
$ws.Range("I6").Value = 229
$ws.Range("I7").Value = 715

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 166
$ws.Range("I7").Value = 325

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 427
$ws.Range("I4").Value = 53
$ws.Range("I7").Value = 1136

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 111

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 214
$ws.Range("I3").Value = 276
$ws.Range("I4").Value = 58
$ws.Range("I6").Value = 399
$ws.Range("I7").Value = 974

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 124
$ws.Range("I7").Value = 338

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 90
$ws.Range("I6").Value = 100
$ws.Range("I7").Value = 304

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 271

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 179
$ws.Range("I6").Value = 222
$ws.Range("I7").Value = 633

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value = 54
$ws.Range("I7").Value = 201

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I5").Value = 3
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I3").Value = 188
$ws.Range("I7").Value = 578

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I2").Value = 45
$ws.Range("I3").Value = 57
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 294

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I3").Value = 82
$ws.Range("I7").Value = 391

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I2").Value = 74
$ws.Range("I7").Value = 231

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 437
$ws.Range("I6").Value = 492
$ws.Range("I7").Value = 1514

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 334

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 62
$ws.Range("I7").Value = 295

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I3").Value = 38
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 148

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 147
$ws.Range("I7").Value = 293

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I6").Value = 219
$ws.Range("I7").Value = 801

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 40
